# "new method to create grimoire"
# A26 held the "has_grimoire" flag; it is replaced by a "zaubern" (casting)
# entry, and the active selection moves down to the next free row (A29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "zaubern"

$ws.Range("A29").Select()
